$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1490.575
$ws.Range("J17").Value = 1490.575
$ws.Range("L17").Value = 4471.725
$ws.Range("N17").Value = -4807.725
$ws.Range("H38").Value = 6712.087
$ws.Range("J38").Value = 8499.75
$ws.Range("L38").Value = 25499.25
$ws.Range("N38").Value = -26243.25
$ws.Range("H40").Value = 5000
$ws.Range("J40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5350
$ws.Range("H43").Value = 6133.143
$ws.Range("J43").Value = 7679.222
$ws.Range("L43").Value = 7679.222
$ws.Range("N43").Value = -7817.222
$ws.Range("H86").Value = 3180.111
$ws.Range("J86").Value = 3156
$ws.Range("L86").Value = 3156
$ws.Range("N86").Value = -5402
$ws.Range("H89").Value = 3180.111
$ws.Range("J89").Value = 3156
$ws.Range("L89").Value = 15780
$ws.Range("N89").Value = -27012
$ws.Range("H116").Value = 4238.5835
$ws.Range("I116").Value = 3651.6667
$ws.Range("K116").Value = 3651.6667
$ws.Range("M116").Value = -209.6667000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2616.3147
$ws.Range("I32").Value = 2500.3333
$ws.Range("J32").Value = 4588
$ws.Range("K32").Value = 2500.3333
$ws.Range("L32").Value = 4588
$ws.Range("M32").Value = -2213.3333
$ws.Range("N32").Value = -5162
$ws.Range("H61").Value = 5126.6113
$ws.Range("I61").Value = 5126.6113
$ws.Range("K61").Value = 5126.6113
$ws.Range("M61").Value = -4914.6113
$ws.Range("H74").Value = 3243.889
$ws.Range("I74").Value = 3399.375
$ws.Range("K74").Value = 3399.375
$ws.Range("M74").Value = -2525.375
$ws.Range("H77").Value = 3243.889
$ws.Range("I77").Value = 3399.375
$ws.Range("K77").Value = 16996.875
$ws.Range("M77").Value = -12628.875
$ws.Range("H132").Value = 1901
$ws.Range("I132").Value = 1889
$ws.Range("K132").Value = 5667
$ws.Range("M132").Value = -3137
$ws.Range("H136").Value = 5126.6113
$ws.Range("I136").Value = 5126.6113
$ws.Range("K136").Value = 15379.8339
$ws.Range("M136").Value = -12829.8339

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 449.4762
$ws.Range("I22").Value = 449.4762
$ws.Range("K22").Value = 449.4762
$ws.Range("M22").Value = -276.4762
$ws.Range("H86").Value = 2427.4666
$ws.Range("I86").Value = 2505
$ws.Range("K86").Value = 2505
$ws.Range("M86").Value = -1382
$ws.Range("H89").Value = 2427.4666
$ws.Range("I89").Value = 2505
$ws.Range("K89").Value = 12525
$ws.Range("M89").Value = -6909
$ws.Range("H132").Value = 74999.5
$ws.Range("J132").Value = 74999.5
$ws.Range("L132").Value = 74999.5
$ws.Range("N132").Value = -85119.5
$ws.Range("H135").Value = 100000
$ws.Range("J135").Value = 100000
$ws.Range("L135").Value = 100000
$ws.Range("N135").Value = -110140
$ws.Range("H140").Value = 75000
$ws.Range("J140").Value = 75000
$ws.Range("L140").Value = 75000
$ws.Range("N140").Value = -85360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1570
$ws.Range("I105").Value = 1500
$ws.Range("J105").Value = 1780
$ws.Range("K105").Value = 1500
$ws.Range("L105").Value = 1780
$ws.Range("M105").Value = 247
$ws.Range("N105").Value = -5274
$ws.Range("H107").Value = 1304.6
$ws.Range("I107").Value = 773.4
$ws.Range("K107").Value = 773.4
$ws.Range("M107").Value = 1146.6
$ws.Range("H132").Value = 2252.027
$ws.Range("I132").Value = 2291.2058
$ws.Range("K132").Value = 6873.617400000001
$ws.Range("M132").Value = -4343.617400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1256.9333
$ws.Range("J113").Value = 1406.8572
$ws.Range("L113").Value = 4220.571599999999
$ws.Range("N113").Value = -8560.571599999999
$ws.Range("H117").Value = 1626.909
$ws.Range("J117").Value = 1324.2858
$ws.Range("L117").Value = 3972.8574
$ws.Range("N117").Value = -10856.8574
$ws.Range("H140").Value = 1363.4117
$ws.Range("I140").Value = 1251.9333
$ws.Range("K140").Value = 3755.7999
$ws.Range("M140").Value = 1424.2001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 48000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 48000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H70").Value = 7750
$ws.Range("I70").Value = 7000
$ws.Range("K70").Value = 7000
$ws.Range("M70").Value = -6730
$ws.Range("H73").Value = 7750
$ws.Range("I73").Value = 7000
$ws.Range("K73").Value = 7000
$ws.Range("M73").Value = -6064
$ws.Range("H80").Value = 6382.3335
$ws.Range("I80").Value = 5530
$ws.Range("J80").Value = 6666.4443
$ws.Range("K80").Value = 5530
$ws.Range("L80").Value = 6666.4443
$ws.Range("M80").Value = -4532
$ws.Range("N80").Value = -8662.444299999999
$ws.Range("H83").Value = 6382.3335
$ws.Range("I83").Value = 5530
$ws.Range("J83").Value = 6666.4443
$ws.Range("K83").Value = 27650
$ws.Range("L83").Value = 33332.2215
$ws.Range("M83").Value = -22658
$ws.Range("N83").Value = -43316.2215
$ws.Range("H102").Value = 4486.2188
$ws.Range("I102").Value = 3984.8965
$ws.Range("K102").Value = 3984.8965
$ws.Range("M102").Value = -2362.8965
$ws.Range("H132").Value = 3842.8572
$ws.Range("I132").Value = 3842.8572
$ws.Range("K132").Value = 11528.5716
$ws.Range("M132").Value = -8998.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4042.75
$ws.Range("I7").Value = 3363.1428
$ws.Range("K7").Value = 3363.1428
$ws.Range("M7").Value = -3251.1428
$ws.Range("H35").Value = 4975
$ws.Range("I35").Value = 1470.8334
$ws.Range("K35").Value = 1470.8334
$ws.Range("M35").Value = -1134.8334
$ws.Range("H63").Value = 48000
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 48000
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H82").Value = 3801.4092
$ws.Range("I82").Value = 3519.9167
$ws.Range("J82").Value = 4139.2
$ws.Range("K82").Value = 3519.9167
$ws.Range("L82").Value = 4139.2
$ws.Range("M82").Value = -3158.9167
$ws.Range("N82").Value = -4861.2
$ws.Range("H85").Value = 3801.4092
$ws.Range("I85").Value = 3519.9167
$ws.Range("J85").Value = 4139.2
$ws.Range("K85").Value = 3519.9167
$ws.Range("L85").Value = 4139.2
$ws.Range("M85").Value = -2271.9167
$ws.Range("N85").Value = -6635.2
$ws.Range("H101").Value = 51340
$ws.Range("J101").Value = 51340
$ws.Range("L101").Value = 51340
$ws.Range("N101").Value = -57830
$ws.Range("H122").Value = 4803.1
$ws.Range("I122").Value = 4792.3335
$ws.Range("J122").Value = 4900
$ws.Range("K122").Value = 14377.0005
$ws.Range("L122").Value = 14700
$ws.Range("M122").Value = -11927.0005
$ws.Range("N122").Value = -19600
$ws.Range("H126").Value = 4042.75
$ws.Range("I126").Value = 3363.1428
$ws.Range("K126").Value = 10089.4284
$ws.Range("M126").Value = -7619.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 110694.5
$ws.Range("J93").Value = 110694.5
$ws.Range("L93").Value = 110694.5
$ws.Range("N93").Value = -115686.5
$ws.Range("H107").Value = 2666.7222
$ws.Range("I107").Value = 2133.4
$ws.Range("K107").Value = 6400.200000000001
$ws.Range("M107").Value = -4480.200000000001
$ws.Range("H122").Value = 4275.3335
$ws.Range("I122").Value = 1918.1177
$ws.Range("K122").Value = 5754.3531
$ws.Range("M122").Value = -3304.3531
$ws.Range("H136").Value = 1322.15
$ws.Range("I136").Value = 1208.8235
$ws.Range("K136").Value = 3626.4705
$ws.Range("M136").Value = -1076.4705

Write-Host "All 205 cell updates applied"
